$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values from 45182 to 45184 for rows 2-10
$ws.Range("C2:C10").Value = 45184
